$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Cells.Item(17, 8).Value = 1954500.8
$ws_ALC.Cells.Item(17, 10).Value = 1954500.8
$ws_ALC.Cells.Item(17, 12).Value = 5863502.4
$ws_ALC.Cells.Item(17, 14).Value = -5863838.4

$ws_ALC.Cells.Item(64, 8).Value = 11749.875
$ws_ALC.Cells.Item(64, 9).Value = 7799.8
$ws_ALC.Cells.Item(64, 10).Value = 18333.334
$ws_ALC.Cells.Item(64, 11).Value = 7799.8
$ws_ALC.Cells.Item(64, 12).Value = 18333.334
$ws_ALC.Cells.Item(64, 13).Value = -7551.8
$ws_ALC.Cells.Item(64, 14).Value = -18829.334

$ws_ALC.Cells.Item(67, 8).Value = 11749.875
$ws_ALC.Cells.Item(67, 9).Value = 7799.8
$ws_ALC.Cells.Item(67, 10).Value = 18333.334
$ws_ALC.Cells.Item(67, 11).Value = 7799.8
$ws_ALC.Cells.Item(67, 12).Value = 18333.334
$ws_ALC.Cells.Item(67, 13).Value = -6941.8
$ws_ALC.Cells.Item(67, 14).Value = -20049.334

$ws_ALC.Cells.Item(100, 8).Value = 1834.2142
$ws_ALC.Cells.Item(100, 9).Value = 1206.0769
$ws_ALC.Cells.Item(100, 10).Value = 10000
$ws_ALC.Cells.Item(100, 11).Value = 1206.0769
$ws_ALC.Cells.Item(100, 12).Value = 10000
$ws_ALC.Cells.Item(100, 13).Value = -665.0769
$ws_ALC.Cells.Item(100, 14).Value = -11082

$ws_ALC.Cells.Item(132, 8).Value = 1438.638
$ws_ALC.Cells.Item(132, 9).Value = 1445.7358
$ws_ALC.Cells.Item(132, 11).Value = 4337.207399999999
$ws_ALC.Cells.Item(132, 13).Value = -1807.207399999999

$ws_ALC.Cells.Item(138, 8).Value = 5159267
$ws_ALC.Cells.Item(138, 9).Value = 2775.5557
$ws_ALC.Cells.Item(138, 10).Value = 5686635.5
$ws_ALC.Cells.Item(138, 11).Value = 8326.667099999999
$ws_ALC.Cells.Item(138, 12).Value = 17059906.5
$ws_ALC.Cells.Item(138, 13).Value = -3186.667099999999
$ws_ALC.Cells.Item(138, 14).Value = -17070186.5

$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Cells.Item(61, 8).Value = 261986.39
$ws_ARM.Cells.Item(61, 9).Value = 2911.65
$ws_ARM.Cells.Item(61, 11).Value = 2911.65
$ws_ARM.Cells.Item(61, 13).Value = -2699.65

$ws_ARM.Cells.Item(63, 8).Value = 2579.8
$ws_ARM.Cells.Item(63, 9).Value = 966.3333
$ws_ARM.Cells.Item(63, 11).Value = 966.3333
$ws_ARM.Cells.Item(63, 13).Value = -280.3333

$ws_ARM.Cells.Item(66, 8).Value = 2579.8
$ws_ARM.Cells.Item(66, 9).Value = 966.3333
$ws_ARM.Cells.Item(66, 11).Value = 4831.6665
$ws_ARM.Cells.Item(66, 13).Value = -1399.6665

$ws_ARM.Cells.Item(102, 8).Value = 3339.625
$ws_ARM.Cells.Item(102, 9).Value = 3529.4285
$ws_ARM.Cells.Item(102, 11).Value = 3529.4285
$ws_ARM.Cells.Item(102, 13).Value = -1907.4285

$ws_ARM.Cells.Item(132, 8).Value = 2619.3948
$ws_ARM.Cells.Item(132, 9).Value = 2304.6333
$ws_ARM.Cells.Item(132, 10).Value = 3799.75
$ws_ARM.Cells.Item(132, 11).Value = 6913.8999
$ws_ARM.Cells.Item(132, 12).Value = 11399.25
$ws_ARM.Cells.Item(132, 13).Value = -4383.8999
$ws_ARM.Cells.Item(132, 14).Value = -16459.25

$ws_ARM.Cells.Item(136, 8).Value = 261986.39
$ws_ARM.Cells.Item(136, 9).Value = 2911.65
$ws_ARM.Cells.Item(136, 11).Value = 8734.950000000001
$ws_ARM.Cells.Item(136, 13).Value = -6184.950000000001

$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Cells.Item(74, 8).Value = 79999
$ws_BSM.Cells.Item(74, 10).Value = 0
$ws_BSM.Cells.Item(74, 12).Value = 0
$ws_BSM.Cells.Item(74, 14).ClearContents()

$ws_BSM.Cells.Item(77, 8).Value = 79999
$ws_BSM.Cells.Item(77, 10).Value = 0
$ws_BSM.Cells.Item(77, 12).Value = 0
$ws_BSM.Cells.Item(77, 14).ClearContents()

$ws_BSM.Cells.Item(81, 8).Value = 33019.668
$ws_BSM.Cells.Item(81, 10).Value = 33019.668
$ws_BSM.Cells.Item(81, 12).Value = 33019.668
$ws_BSM.Cells.Item(81, 14).Value = -35141.668

$ws_BSM.Cells.Item(84, 8).Value = 33019.668
$ws_BSM.Cells.Item(84, 10).Value = 33019.668
$ws_BSM.Cells.Item(84, 12).Value = 99059.00399999999
$ws_BSM.Cells.Item(84, 14).Value = -109667.004

$ws_BSM.Cells.Item(99, 8).Value = 7145
$ws_BSM.Cells.Item(99, 9).Value = 4340.7144
$ws_BSM.Cells.Item(99, 11).Value = 4340.7144
$ws_BSM.Cells.Item(99, 13).Value = -2842.7144

$ws_BSM.Cells.Item(104, 8).Value = 0
$ws_BSM.Cells.Item(104, 10).Value = 0
$ws_BSM.Cells.Item(104, 12).Value = 0
$ws_BSM.Cells.Item(104, 14).ClearContents()

$ws_BSM.Cells.Item(105, 8).Value = 2726.8948
$ws_BSM.Cells.Item(105, 9).Value = 2825.6875
$ws_BSM.Cells.Item(105, 10).Value = 2200
$ws_BSM.Cells.Item(105, 11).Value = 2825.6875
$ws_BSM.Cells.Item(105, 12).Value = 2200
$ws_BSM.Cells.Item(105, 13).Value = -1078.6875
$ws_BSM.Cells.Item(105, 14).Value = -5694

$ws_BSM.Cells.Item(107, 8).Value = 1608.7838
$ws_BSM.Cells.Item(107, 9).Value = 1671.3529
$ws_BSM.Cells.Item(107, 11).Value = 1671.3529
$ws_BSM.Cells.Item(107, 13).Value = 248.6470999999999

$ws_BSM.Cells.Item(134, 8).Value = 1225.175
$ws_BSM.Cells.Item(134, 9).Value = 1205.3077
$ws_BSM.Cells.Item(134, 11).Value = 3615.9231
$ws_BSM.Cells.Item(134, 13).Value = -1080.9231

$ws_BSM.Cells.Item(138, 8).Value = 102449.8
$ws_BSM.Cells.Item(138, 10).Value = 102449.8
$ws_BSM.Cells.Item(138, 12).Value = 102449.8
$ws_BSM.Cells.Item(138, 14).Value = -112729.8

$ws_BSM.Cells.Item(139, 8).Value = 100000
$ws_BSM.Cells.Item(139, 10).Value = 100000
$ws_BSM.Cells.Item(139, 12).Value = 100000
$ws_BSM.Cells.Item(139, 14).Value = -110280

$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Cells.Item(7, 8).Value = 760.6842
$ws_CRP.Cells.Item(7, 10).Value = 467.625
$ws_CRP.Cells.Item(7, 12).Value = 467.625
$ws_CRP.Cells.Item(7, 14).Value = -693.625

$ws_CRP.Cells.Item(16, 8).Value = 7865.3335
$ws_CRP.Cells.Item(16, 9).Value = 5998.5386
$ws_CRP.Cells.Item(16, 11).Value = 5998.5386
$ws_CRP.Cells.Item(16, 13).Value = -5711.5386

$ws_CRP.Cells.Item(31, 8).Value = 86637.664
$ws_CRP.Cells.Item(31, 9).Value = 102864.1
$ws_CRP.Cells.Item(31, 10).Value = 5505.5
$ws_CRP.Cells.Item(31, 11).Value = 102864.1
$ws_CRP.Cells.Item(31, 12).Value = 5505.5
$ws_CRP.Cells.Item(31, 13).Value = -102569.1
$ws_CRP.Cells.Item(31, 14).Value = -6095.5

$ws_CRP.Cells.Item(34, 8).Value = 86637.664
$ws_CRP.Cells.Item(34, 9).Value = 102864.1
$ws_CRP.Cells.Item(34, 10).Value = 5505.5
$ws_CRP.Cells.Item(34, 11).Value = 102864.1
$ws_CRP.Cells.Item(34, 12).Value = 5505.5
$ws_CRP.Cells.Item(34, 13).Value = -102662.1
$ws_CRP.Cells.Item(34, 14).Value = -5909.5

$ws_CRP.Cells.Item(113, 8).Value = 7865.3335
$ws_CRP.Cells.Item(113, 9).Value = 5998.5386
$ws_CRP.Cells.Item(113, 11).Value = 5998.5386
$ws_CRP.Cells.Item(113, 13).Value = -3828.5386

$ws_CRP.Cells.Item(134, 8).Value = 9259.394
$ws_CRP.Cells.Item(134, 9).Value = 5094.926
$ws_CRP.Cells.Item(134, 10).Value = 27999.5
$ws_CRP.Cells.Item(134, 11).Value = 15284.778
$ws_CRP.Cells.Item(134, 12).Value = 83998.5
$ws_CRP.Cells.Item(134, 13).Value = -12749.778
$ws_CRP.Cells.Item(134, 14).Value = -89068.5

$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Cells.Item(7, 8).Value = 97.882355
$ws_CUL.Cells.Item(7, 9).Value = 105.53333
$ws_CUL.Cells.Item(7, 10).Value = 40.5
$ws_CUL.Cells.Item(7, 11).Value = 316.59999
$ws_CUL.Cells.Item(7, 12).Value = 121.5
$ws_CUL.Cells.Item(7, 13).Value = -204.59999
$ws_CUL.Cells.Item(7, 14).Value = -345.5

$ws_CUL.Cells.Item(96, 8).Value = 14666.667
$ws_CUL.Cells.Item(96, 10).Value = 14666.667
$ws_CUL.Cells.Item(96, 12).Value = 44000.001
$ws_CUL.Cells.Item(96, 14).Value = -48118.001

$ws_CUL.Cells.Item(113, 8).Value = 1287.4615
$ws_CUL.Cells.Item(113, 10).Value = 1273.4445
$ws_CUL.Cells.Item(113, 12).Value = 3820.3335
$ws_CUL.Cells.Item(113, 14).Value = -8160.333500000001

$ws_CUL.Cells.Item(122, 8).Value = 2838.7222
$ws_CUL.Cells.Item(122, 9).Value = 3248.8333
$ws_CUL.Cells.Item(122, 10).Value = 2633.6667
$ws_CUL.Cells.Item(122, 11).Value = 29239.4997
$ws_CUL.Cells.Item(122, 12).Value = 23703.0003
$ws_CUL.Cells.Item(122, 13).Value = -26789.4997
$ws_CUL.Cells.Item(122, 14).Value = -28603.0003

$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Cells.Item(102, 8).Value = 34484344
$ws_GSM.Cells.Item(102, 9).Value = 1645.2963
$ws_GSM.Cells.Item(102, 11).Value = 1645.2963
$ws_GSM.Cells.Item(102, 13).Value = -23.29629999999997

$ws_GSM.Cells.Item(111, 8).Value = 0
$ws_GSM.Cells.Item(111, 10).Value = 0
$ws_GSM.Cells.Item(111, 12).Value = 0
$ws_GSM.Cells.Item(111, 14).ClearContents()

$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Cells.Item(87, 8).Value = 199000
$ws_LTW.Cells.Item(87, 10).Value = 199000
$ws_LTW.Cells.Item(87, 12).Value = 199000
$ws_LTW.Cells.Item(87, 14).Value = -201246

$ws_LTW.Cells.Item(90, 8).Value = 199000
$ws_LTW.Cells.Item(90, 10).Value = 199000
$ws_LTW.Cells.Item(90, 12).Value = 597000
$ws_LTW.Cells.Item(90, 14).Value = -608232

$ws_LTW.Cells.Item(136, 8).Value = 5113.067
$ws_LTW.Cells.Item(136, 9).Value = 4259.2
$ws_LTW.Cells.Item(136, 10).Value = 6820.8
$ws_LTW.Cells.Item(136, 11).Value = 12777.6
$ws_LTW.Cells.Item(136, 12).Value = 20462.4
$ws_LTW.Cells.Item(136, 13).Value = -10227.6
$ws_LTW.Cells.Item(136, 14).Value = -25562.4

$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Cells.Item(33, 8).Value = 20000
$ws_WVR.Cells.Item(33, 9).Value = 0
$ws_WVR.Cells.Item(33, 11).Value = 0
$ws_WVR.Cells.Item(33, 13).ClearContents()

$ws_WVR.Cells.Item(36, 8).Value = 20000
$ws_WVR.Cells.Item(36, 9).Value = 0
$ws_WVR.Cells.Item(36, 11).Value = 0
$ws_WVR.Cells.Item(36, 13).ClearContents()

$ws_WVR.Cells.Item(62, 8).Value = 12329.111
$ws_WVR.Cells.Item(62, 9).Value = 8368.25
$ws_WVR.Cells.Item(62, 11).Value = 8368.25
$ws_WVR.Cells.Item(62, 13).Value = -7744.25

$ws_WVR.Cells.Item(65, 8).Value = 12329.111
$ws_WVR.Cells.Item(65, 9).Value = 8368.25
$ws_WVR.Cells.Item(65, 11).Value = 41841.25
$ws_WVR.Cells.Item(65, 13).Value = -38721.25

$ws_WVR.Cells.Item(96, 8).Value = 1968.625
$ws_WVR.Cells.Item(96, 9).Value = 1869.8
$ws_WVR.Cells.Item(96, 10).Value = 2133.3333
$ws_WVR.Cells.Item(96, 11).Value = 1869.8
$ws_WVR.Cells.Item(96, 12).Value = 2133.3333
$ws_WVR.Cells.Item(96, 13).Value = -496.8
$ws_WVR.Cells.Item(96, 14).Value = -4879.3333

$ws_WVR.Cells.Item(126, 8).Value = 3904.7222
$ws_WVR.Cells.Item(126, 9).Value = 3726.6155
$ws_WVR.Cells.Item(126, 10).Value = 4367.8
$ws_WVR.Cells.Item(126, 11).Value = 11179.8465
$ws_WVR.Cells.Item(126, 12).Value = 13103.4
$ws_WVR.Cells.Item(126, 13).Value = -8709.8465
$ws_WVR.Cells.Item(126, 14).Value = -18043.4

$ws_WVR.Cells.Item(130, 8).Value = 18703
$ws_WVR.Cells.Item(130, 10).Value = 18703
$ws_WVR.Cells.Item(130, 12).Value = 18703
$ws_WVR.Cells.Item(130, 14).Value = -28743

$ws_WVR.Cells.Item(132, 8).Value = 1877.3334
$ws_WVR.Cells.Item(132, 9).Value = 1544.2142
$ws_WVR.Cells.Item(132, 10).Value = 4542.2856
$ws_WVR.Cells.Item(132, 11).Value = 4632.642599999999
$ws_WVR.Cells.Item(132, 12).Value = 13626.8568
$ws_WVR.Cells.Item(132, 13).Value = -2102.642599999999
$ws_WVR.Cells.Item(132, 14).Value = -18686.8568
